$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.596.49"
$ws.Range("E2").Value = "  +0.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.467.18"
$ws.Range("E3").Value = "  +0.13%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.82"
$ws.Range("E5").Value = "  +1.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.30"
$ws.Range("E6").Value = "  +1.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  +0.60%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  +0.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0868"
$ws.Range("E10").Value = "  +9.41%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.95"
$ws.Range("E11").Value = "  +1.28%  "

# Row 12
$ws.Range("E12").Value = "  +0.09%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.846.94"
$ws.Range("E13").Value = "  +0.11%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("E14").Value = "  +0.02%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.46"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.489.00"
$ws.Range("E16").Value = "  +1.55%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").Value = "  +2.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.548.13"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("E19").Value = "  -0.81%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0944"
$ws.Range("E20").Value = "  +0.59%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.84"
$ws.Range("E21").Value = "  -0.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.26"
$ws.Range("E22").Value = "  -0.47%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.98"
$ws.Range("E23").Value = "  +0.97%  "

# Row 24
$ws.Range("E24").Value = "  +1.28%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.95"
$ws.Range("E25").Value = "  +2.34%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.75"
$ws.Range("E27").Value = "  +2.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("E28").Value = "  -0.81%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  +0.47%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.35"
$ws.Range("E30").Value = "  +3.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.32"
$ws.Range("E31").Value = "  +0.80%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.46"
$ws.Range("E32").Value = "  +0.52%  "

# Row 33
$ws.Range("E33").Value = "  -0.10%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0762"
$ws.Range("E34").Value = "  +0.65%  "

# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.56"
$ws.Range("E35").Value = "  -0.49%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.23"
$ws.Range("E36").Value = "  +0.12%  "

# Row 37
$ws.Range("E37").Value = "  +4.19%  "

# Row 38
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.116"
$ws.Range("E38").Value = "  +1.61%  "

# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.89"
$ws.Range("E39").Value = "  +1.07%  "

# Row 40
$ws.Range("E40").Value = "  +1.12%  "

# Row 41
$ws.Range("E41").Value = "  +7.62%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.00"
$ws.Range("E42").Value = "  +0.16%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.988.64"
$ws.Range("E43").Value = "  +1.60%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0283"
$ws.Range("E44").Value = "  +0.41%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.79"
$ws.Range("E45").Value = "  +0.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.95"
$ws.Range("E46").Value = "  +1.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.49"
$ws.Range("E47").Value = "  +5.45%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.705.00"
$ws.Range("E48").Value = "  +0.09%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.51"
$ws.Range("E49").Value = "  +0.81%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.48"
$ws.Range("E50").Value = "  +5.49%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.66"
$ws.Range("E51").Value = "  -0.43%  "
